$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old header row and the now-unused K:L columns entirely.
$ws.Range("A1:L1").ClearContents()
$ws.Range("K1:L10").ClearContents()

# New single header cell.
$ws.Range("A1").Value = "Herald College Kathmandu"

# Row data: Day, Time, Hours, Module Code, Module Title, Class Type,
# Lecturer, Group, Block, Room  (columns A..J)
$rows = @(
    @{ A="SUN"; B="12:00-14:30"; C=2.5; D="5CS024"; E="Collaborative Development"; F="Workshop"; G="Mr. Biraj Dulal"; H="L5CG7"; I="WLV"; J="TR-01 Dudley" },
    @{ A="MON"; B="7:00-9:30";   C=2.5; D="5CS022"; E="Human Computer Interaction"; F="Workshop"; G="Mr. Ayush Shakya"; H="L5CG7"; I="HCK"; J="Lab-05 Basantapur" },
    @{ A="TUE"; B="7:00-9:00";   C=2;   D="5CS022"; E="Human Computer Interaction"; F="Lecture";  G="Mr. Apurba Neupane"; H="L5CG(5+6+7+8)"; I="WLV"; J="LT-02 Telford" },
    @{ A="TUE"; B="9:30-11:30";  C=2;   D="5CS020"; E="Distributed and Cloud Systems Programming"; F="Lecture"; G="Mr. Sumanta Silwal"; H="L5CG(5+6+7+8)"; I="WLV"; J="LT-01 Wulfruna" },
    @{ A="WED"; B="7:00-9:00";   C=2;   D="5CS024"; E="Collaborative Development"; F="Lecture"; G="Mr. Raj Shrestha"; H="L5CG(5+6+7+8)"; I="WLV"; J="LT-02 Telford" },
    @{ A="WED"; B="12:00-14:00"; C=2;   D="5CS020"; E="Distributed and Cloud Systems Programming"; F="Tutorial"; G="Mr. Shishir Poudel"; H="L5CG7"; I="WLV"; J="TR-01 Dudley" },
    @{ A="THU"; B="7:00-9:00";   C=2;   D="5CS022"; E="Human Computer Interaction"; F="Tutorial"; G="Mr. Ayush Shakya"; H="L5CG7"; I="WLV"; J="TR-01 Dudley" },
    @{ A="FRI"; B="9:00-11:00";  C=2;   D="5CS024"; E="Collaborative Development"; F="Tutorial"; G="Mr. Biraj Dulal"; H="L5CG7"; I="WLV"; J="TR-03 Westbromwich" },
    @{ A="FRI"; B="12:30-15:30"; C=2.5; D="5CS020"; E="Distributed and Cloud Systems Programming"; F="Workshop"; G="Mr. Shishir Poudel"; H="L5CG7"; I="WLV"; J="Lab-02 Moseley" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $r++
}
